# Refresh the cryptocurrency Price (D) / Volume(1h) (E) columns on Sheet1,
# matching the automated "Updated symbol list" GitHub Actions data refresh.
# The source data cells are plain text (e.g. "329.07", "1.71%"), not numbers,
# so each new value must be written back as text too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    'D2' = '329.07';
    'E2' = '1.71%';
    'D3' = '41.49';
    'E3' = '4.65%';
    'D4' = '5.624';
    'E4' = '-4.36%';
    'D5' = '0.08170';
    'E5' = '1.72%';
    'D6' = '2.026';
    'E6' = '4.11%';
    'E7' = '1.11%';
    'E8' = '-0.97%';
    'D9' = '2.942';
    'E9' = '-0.19%';
    'D10' = '0.9180';
    'E10' = '-1.23%';
    'D11' = '0.1271';
    'E11' = '-0.21%';
    'D12' = '0.1950';
    'E12' = '-0.78%';
    'D13' = '0.09280';
    'E13' = '1.70%';
    'D14' = '0.03740';
    'E14' = '5.29%';
    'D15' = '0.1061';
    'E15' = '1.41%';
    'D16' = '0.001308';
    'E16' = '1.36%';
    'D17' = '0.006220';
    'E17' = '-0.30%';
    'D19' = '3.439';
    'E19' = '2.71%';
    'E20' = '-1.20%';
    'D21' = '8.271';
    'E21' = '-5.09%';
    'D22' = '0.1394';
    'E22' = '1.65%';
    'E23' = '-2.36%';
    'D24' = '0.04427';
    'E24' = '0.37%';
    'D25' = '0.001261';
    'E25' = '-0.20%';
    'D26' = '0.004308';
    'E26' = '-1.63%';
    'D27' = '0.0001181';
    'E27' = '3.58%';
    'D39' = '0.02764';
    'E39' = '9.31%';
    'D40' = '0.05414';
    'E40' = '2.86%';
    'D41' = '0.007673';
    'E41' = '3.12%';
    'D42' = '0.1414';
    'E42' = '0.57%';
    'D43' = '0.008999';
    'E43' = '-6.36%';
    'D44' = '0.002122';
    'E44' = '0.23%';
    'D45' = '0.01143';
    'E45' = '14.57%';
    'D46' = '0.00006790';
    'E46' = '1.01%';
    'E47' = '0.14%';
    'D48' = '0.002283';
    'E48' = '-0.34%';
    'D49' = '0.003500';
    'E49' = '16.59%';
    'D50' = '0.00002104';
    'E50' = '0.14%';
    'E51' = '0.14%'
}

foreach ($addr in $newValues.Keys) {
    $cell = $ws.Range($addr)
    # A leading apostrophe forces Excel to store the value as literal text
    # (otherwise "329.07" / "1.71%" would be auto-converted to a number /
    # percentage), matching the inline-string cells already in the sheet.
    $cell.Value = "'" + $newValues[$addr]
    # Writing a quote-prefixed value tags the cell with a quote-prefix style;
    # put the style back to the workbook default so no formatting changes
    # beyond the text itself are introduced.
    $cell.Style = "Normal"
}
